# Water Quality Data test-parser workbook: purge more pytz calls, fix test parser header
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Time (24HR)" column header is being simplified to just "Time"
$ws.Range("D2").Value = "Time"

# Leave the active cell/selection on D5, matching the saved workbook view
$ws.Range("D5").Select() | Out-Null
